$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2.32
$ws.Range("I2").Value = 2.38
$ws.Range("J2").Value = 3.85
$ws.Range("N2").Value = 4.8
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 2.26
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.5
$ws.Range("S2").Value = 2.98
$ws.Range("T2").Value = 1.65
$ws.Range("U2").Value = 2.46
$ws.Range("X2").Value = 19.5
$ws.Range("Z2").Value = 16
$ws.Range("AB2").Value = 15
$ws.Range("AC2").Value = 8.6
$ws.Range("AE2").Value = 23
$ws.Range("AF2").Value = 23
$ws.Range("AG2").Value = 13
$ws.Range("AI2").Value = 32
$ws.Range("AJ2").Value = 50
$ws.Range("AN2").Value = 24
$ws.Range("AO2").Value = 15.5
# Row 3
$ws.Range("F3").Value = 3.5
$ws.Range("G3").Value = 4.1
$ws.Range("H3").Value = 2.08
$ws.Range("I3").Value = 2.24
$ws.Range("J3").Value = 3.45
$ws.Range("K3").Value = 3.9
$ws.Range("L3").Value = 1.45
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 3.6
$ws.Range("O3").Value = 1.37
$ws.Range("P3").Value = 1.87
$ws.Range("R3").Value = 1.35
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 1.33
$ws.Range("X3").Value = 17.5
$ws.Range("Z3").Value = 1000
$ws.Range("AH3").Value = 980
# Row 4
$ws.Range("G4").Value = 2.94
$ws.Range("H4").Value = 2.66
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 3.6
$ws.Range("L4").Value = 1.41
$ws.Range("N4").Value = 3.9
$ws.Range("Q4").Value = 1.92
$ws.Range("R4").Value = 1.43
$ws.Range("T4").Value = 1.69
$ws.Range("V4").Value = 1.55
$ws.Range("W4").Value = 1.5
$ws.Range("X4").Value = 1000
$ws.Range("AA4").Value = 60
$ws.Range("AD4").Value = 990
$ws.Range("AF4").Value = 1000
$ws.Range("AH4").Value = 19
$ws.Range("AN4").Value = 32
# Row 5
$ws.Range("I5").Value = 16
$ws.Range("K5").Value = 5.1
$ws.Range("L5").Value = 1.45
$ws.Range("O5").Value = 1.38
$ws.Range("P5").Value = 1.75
$ws.Range("Q5").Value = 2.14
$ws.Range("S5").Value = 4
$ws.Range("T5").Value = 2.72
$ws.Range("U5").Value = 1.47
$ws.Range("AB5").Value = 6
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 75
$ws.Range("AF5").Value = 6.6
$ws.Range("AJ5").Value = 10.5
$ws.Range("AL5").Value = 90
$ws.Range("AN5").Value = 8.800000000000001
# Row 6
$ws.Range("I6").Value = 4.6
$ws.Range("O6").Value = 1.28
$ws.Range("S6").Value = 3.05
$ws.Range("T6").Value = 1.76
$ws.Range("U6").Value = 2.18
$ws.Range("W6").Value = 2.02
$ws.Range("X6").Value = 20
$ws.Range("Z6").Value = 95
$ws.Range("AE6").Value = 200
$ws.Range("AG6").Value = 12
$ws.Range("AH6").Value = 19
$ws.Range("AJ6").Value = 26
$ws.Range("AK6").Value = 23
$ws.Range("AL6").Value = 38
$ws.Range("AN6").Value = 14.5
$ws.Range("AO6").Value = 600
# Row 7
$ws.Range("F7").Value = 1.43
$ws.Range("I7").Value = 8.6
$ws.Range("J7").Value = 5.1
$ws.Range("M7").Value = 1.03
$ws.Range("Q7").Value = 1.48
$ws.Range("R7").Value = 1.75
$ws.Range("U7").Value = 2.34
$ws.Range("AA7").Value = 230
$ws.Range("AC7").Value = 14
$ws.Range("AD7").Value = 80
# Row 8
$ws.Range("F8").Value = 1.66
$ws.Range("G8").Value = 1.7
$ws.Range("I8").Value = 8.6
$ws.Range("J8").Value = 3.4
$ws.Range("K8").Value = 3.8
$ws.Range("N8").Value = 2.44
$ws.Range("P8").Value = 1.48
$ws.Range("Q8").Value = 2.88
$ws.Range("U8").Value = 1.52
$ws.Range("W8").Value = 2.38
$ws.Range("AN8").Value = 24
# Row 9
$ws.Range("F9").Value = 1.78
$ws.Range("G9").Value = 1.88
$ws.Range("J9").Value = 3.7
$ws.Range("K9").Value = 3.9
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 3.45
$ws.Range("P9").Value = 1.8
$ws.Range("R9").Value = 1.3
$ws.Range("T9").Value = 1.92
$ws.Range("V9").Value = 1.21
$ws.Range("AB9").Value = 14.5
# Row 10
$ws.Range("F10").Value = 1.64
$ws.Range("G10").Value = 1.73
$ws.Range("H10").Value = 5.9
$ws.Range("I10").Value = 7.8
$ws.Range("J10").Value = 3.65
$ws.Range("K10").Value = 4.3
$ws.Range("N10").Value = 3.7
$ws.Range("O10").Value = 1.32
$ws.Range("P10").Value = 1.9
$ws.Range("Q10").Value = 1.96
$ws.Range("S10").Value = 3.5
$ws.Range("T10").Value = 1.89
$ws.Range("V10").Value = 1.16
$ws.Range("W10").Value = 2.36
$ws.Range("Y10").Value = 1000
$ws.Range("AC10").Value = 13
$ws.Range("AF10").Value = 40
# Row 11
$ws.Range("G11").Value = 2.32
$ws.Range("N11").Value = 3.4
$ws.Range("O11").Value = 1.34
$ws.Range("P11").Value = 1.81
$ws.Range("Q11").Value = 2.04
$ws.Range("R11").Value = 1.3
$ws.Range("S11").Value = 3.6
$ws.Range("T11").Value = 1.9
$ws.Range("AF11").Value = 15
# Row 12
$ws.Range("F12").Value = 2.42
$ws.Range("H12").Value = 3.4
$ws.Range("K12").Value = 3.3
$ws.Range("L12").Value = 1.53
$ws.Range("M12").Value = 1.11
$ws.Range("O12").Value = 1.45
$ws.Range("T12").Value = 1.94
$ws.Range("U12").Value = 1.9
$ws.Range("V12").Value = 1.4
$ws.Range("W12").Value = 1.67
$ws.Range("X12").Value = 9.800000000000001
$ws.Range("Y12").Value = 14.5
$ws.Range("Z12").Value = 25
$ws.Range("AC12").Value = 7.6
$ws.Range("AD12").Value = 15
$ws.Range("AE12").Value = 50
$ws.Range("AF12").Value = 14.5
$ws.Range("AG12").Value = 12
$ws.Range("AH12").Value = 22
$ws.Range("AI12").Value = 65
$ws.Range("AJ12").Value = 34
$ws.Range("AK12").Value = 34
$ws.Range("AN12").Value = 30
$ws.Range("AO12").Value = 60
# Row 13
$ws.Range("F13").Value = 2.02
$ws.Range("S13").Value = 3.95
$ws.Range("T13").Value = 1.84
$ws.Range("U13").Value = 1.96
$ws.Range("AB13").Value = 15.5
$ws.Range("AD13").Value = 46
$ws.Range("AJ13").Value = 900
